$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows (health, police) with new values
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 26
$ws.Range("D2").Value = 7.692307692307693

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 0

# Insert a new "fire" row at row 4, pushing "schools" down to row 5
$ws.Range("A4:D4").Insert()

$ws.Range("A4").Value = "fire"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 0

# Update the (now shifted) schools row, which is row 5
$ws.Range("A5").Value = "schools"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 138
$ws.Range("D5").Value = 2.173913043478261
